$d = $word.ActiveDocument

# Apply edits from the end of the document towards the start so that
# earlier (lower) character offsets stay valid as later edits change
# the length of the text that follows them.

$d.Range(1166, 1210).Text = "            Console.WriteLine(`"No iteration "
$d.Range(1081, 1118).Text = "                Console.WriteLine(i);"
$d.Range(1021, 1066).Text = "            for (int i = number; i <= 0; i++)"
$d.Range(920, 957).Text   = "                Console.WriteLine(i);"
$d.Range(860, 905).Text   = "            for (int i = number; i >= 0; i--)"
$d.Range(757, 825).Text   = "            Console.WriteLine(`"Iterating from the number to zero:`");"
$d.Range(693, 746).Text   = "            Console.WriteLine(`"The number is zero.`");"
$d.Range(602, 659).Text   = "            Console.WriteLine(`"The number is negative.`");"
$d.Range(495, 552).Text   = "            Console.WriteLine(`"The number is positive.`");"
$d.Range(439, 460).Text   = "(Console.ReadLine());"
$d.Range(373, 408).Text   = " Console.Write(`"Enter a number: `");"
$d.Range(353, 360).Text   = "] args)"
$d.Range(45, 52).Text     = " : Task 3"
$d.Range(42, 45).Text     = "2 C#"
